# "am facut 4 filtre pe imagine"
# Bumped the "săpt. 4" (week 4 / column F) attendance count for a batch of
# students, renamed one student, and moved the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename student (shared string) Alina Schmidt -> Anita Schmidt
$ws.Range("B46").Value = "Anita Schmidt"

# Column F ("săpt. 4") updates - set new attendance counts.
# The Q column ("Prezențe") holds a shared SUM formula over the row and
# recalculates automatically once F changes.
$f4Updates = @{
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    11 = 2
    12 = 2
    13 = 1
    17 = 1
    19 = 1
    20 = 2
    24 = 2
    25 = 1
    29 = 2
    33 = 1
    34 = 2
    35 = 1
    37 = 2
    38 = 2
    46 = 2
}

foreach ($row in $f4Updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $f4Updates[$row]
}

# Move the active selection to H19 (previously H46)
$ws.Range("H19").Select()
